# Update cryptocurrency price/volume data in the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.936.56"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.818.49"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'309.79"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "'0.3704"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").Value = "'0.07387"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'0.8737"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "'20.50"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.932.32"
$ws.Range("E12").Value = "  +7.10%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'93.19"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.373"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'0.07072"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "'6.516"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'0.000008737"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'14.78"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "26.988.00"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'5.333"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "2.044.81"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'151.67"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'2.225"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "'18.44"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'5.339"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("D30").Value = "'115.72"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "'0.7713"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'1.170"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'4.501"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'2.909"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'1.0000"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "'0.01965"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'0.05291"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.326"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.944"
$ws.Range("D42").Value = "'0.5361"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "'2.373"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'0.1674"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").Value = "'8.472"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'10.49"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "'1.677"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'102.98"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'0.06297"
$ws.Range("E51").Value = "  -0.33%  "
